$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final alphabetically sorted list of card names that go below the title row (A1)
$cards = @(
    "Bounty of Skemfar",
    "Cosmic Intervention",
    "Crown of Skemfar",
    "Elderfang Venom",
    "Ethereal Valkyrie",
    "Hero of Bretagard",
    "Inspired Sphinx",
    "Lathril, Blade of the Elves",
    "Pact of the Serpent",
    "Ranar the Ever-Watchful",
    "Ruthless Winnower",
    "Sage of the Beyond",
    "Serpent's Soul-Jar",
    "Spectral Deluge",
    "Stoic Farmer",
    "Tales of the Ancestors",
    "Wolverine Riders"
)

for ($i = 0; $i -lt $cards.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $cards[$i]
}
